# agregue set de validacion
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Update column D values on Sheet1 ---
$sheet1Updates = @{
    2  = 0.1225
    3  = 0.1225
    4  = 0.1225
    5  = 0.1225
    6  = 0.1225
    7  = 0.1225
    8  = 0.1225
    9  = 0.1225
    10 = 0.2675
    11 = 0.2675
    12 = 0.2675
    13 = 0.2675
    14 = 0.2675
    15 = 0.2675
    16 = 0.2675
    17 = 0.2675
    18 = 0.195
    19 = 0.195
    20 = 0.195
    21 = 0.34
    22 = 0.195
    23 = 0.195
    24 = 0.195
    25 = 0.195
    27 = 0.195
    28 = 0.195
    29 = 0.195
    30 = 0.195
    31 = 0.195
    32 = 0.195
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 4).Value = $sheet1Updates[$row]
}

# --- Add a new Sheet2 with a validation dataset, right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Headers (same column labels as Sheet1)
$ws2.Cells.Item(1, 1).Value = "Cr ppm"
$ws2.Cells.Item(1, 2).Value = "Ni ppm"
$ws2.Cells.Item(1, 3).Value = "Fe ppm"
$ws2.Cells.Item(1, 4).Value = "Co ppm"
$ws2.Cells.Item(1, 5).Value = "B ppm"
$ws2.Range("A1:E1").Style = $ws1.Range("A1:E1").Style

$sheet2Data = @(
    @(143.75, 336.25, 307.5, 0.068125, 0.05249999999999999),
    @(181.25, 238.75, 482.5, 0.104375, 0.0325),
    @(193.75, 303.75, 342.5, 0.249375, 0.0375),
    @(131.25, 206.25, 377.5, 0.285625, 0.0475),
    @(106.25, 368.75, 412.5, 0.321875, 0.0225),
    @(156.25, 433.75, 447.5, 0.213125, 0.0575),
    @(168.75, 401.25, 272.5, 0.176875, 0.0275),
    @(118.75, 271.25, 237.5, 0.140625, 0.0425)
)

$r = 2
foreach ($row in $sheet2Data) {
    for ($c = 0; $c -lt 5; $c++) {
        $ws2.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

$ws1.Select()
